# Refresh Neg_Change and Pos_Change sheets with the latest filtered market snapshot.
$wb = $excel.ActiveWorkbook

$sheet1Data = @(
    @("TATAMOTORS", 721.15, 722.25, 713.5, 715.2, 9409853, 20531027, -0.5416764587567879, "TATAMOTORS"),
    @("UNITDSPR", 1300, 1307.2, 1291.7, 1305, 411404, 822646, -0.499901537234728, "UNITDSPR"),
    @("MOTHERSON", 99.1, 99.5, 97.25, 98.8, 22643077, 56299081, -0.5978073425390372, "MOTHERSON"),
    @("DIVISLAB", 6025, 6050, 5997.5, 6005, 254569, 519368, -0.5098485081868733, "DIVISLAB"),
    @("TVSMOTOR", 3591.3, 3597.3, 3544.4, 3574, 675896, 1406946, -0.5196006101158112, "TVSMOTOR"),
    @("JINDALSTEL", 1043, 1044.45, 1023.2, 1033, 766328, 1823229, -0.5796863696222472, "JINDALSTEL"),
    @("ABFRL", 91.9, 92.48, 89.5, 91.05, 18714097, 40149600, -0.5338908233207803, "ABFRL"),
    @("PRESTIGE", 1539.9, 1550.9, 1530, 1545, 206951, 513985, -0.5973598451316673, "PRESTIGE"),
    @("ATGL", 595, 598.2, 592.05, 594.5, 217589, 431824, -0.4961164733780429, "ATGL"),
    @("DALBHARAT", 2421.8, 2421.8, 2403.1, 2410.5, 142687, 315525, -0.5477790983281832, "DALBHARAT"),
    @("AUBANK", 698, 699.75, 694.85, 696.55, 898505, 1949781, -0.539176451098867, "AUBANK"),
    @("ASHOKLEY", 137.54, 137.7, 134.82, 135.06, 15017284, 33590663, -0.5529327896862292, "ASHOKLEY"),
    @("PNBHOUSING", 803, 810, 799.2, 800.6, 1103494, 2561729, -0.5692385884689598, "PNBHOUSING"),
    @("ANGELONE", 2350, 2356, 2323.1, 2336.6, 338579, 786620, -0.5695774325595586, "ANGELONE")
)

$sheet2Data = @(
    @("BEL", 371.9, 374, 369.9, 373, 10598912, 6897641, 0.5365995417853727, "BEL"),
    @("LT", 3547.2, 3550, 3516.4, 3527.5, 1522116, 1039999, 0.4635744842062348, "LT"),
    @("ZYDUSLIFE", 1010, 1038.7, 1009.2, 1038.5, 1228793, 777187, 0.5810776556993362, "ZYDUSLIFE"),
    @("ABB", 5111, 5150, 5091, 5139.9, 150868, 101582, 0.485184383059991, "ABB"),
    @("GAIL", 172.39, 173.64, 171.7, 173.4, 8186535, 5686270, 0.4397021245913402, "GAIL"),
    @("HAL", 4450, 4474, 4427.7, 4454, 888000, 599045, 0.4823594220801442, "HAL"),
    @("GLENMARK", 2073, 2119.9, 2062, 2115.7, 1269358, 842874, 0.505987846344768, "GLENMARK"),
    @("LUPIN", 1937.1, 1964.5, 1934.7, 1961.9, 486792, 329889, 0.4756236188536144, "LUPIN"),
    @("POONAWALLA", 430, 435, 426.6, 430, 424106, 269428, 0.5740977181287765, "POONAWALLA"),
    @("NUVAMA", 6451.5, 6481, 6346, 6377.5, 83920, 58511, 0.4342602245731572, "NUVAMA")
)

$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# Neg_Change previously had 18 data rows (rows 2-19); the refreshed snapshot only has 14 (rows 2-15).
# Drop the now-unused trailing rows first so the sheet dimension shrinks to A1:I15.
$ws1.Range("A16:I19").EntireRow.Delete() | Out-Null

# Pos_Change previously had 11 data rows (rows 2-12); the refreshed snapshot only has 10 (rows 2-11).
$ws2.Range("A12:I12").EntireRow.Delete() | Out-Null

# Write the refreshed rows for Neg_Change starting at row 2.
$r = 2
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $ws1.Cells.Item($r, 9).Value = $row[8]
    $r++
}

# Write the refreshed rows for Pos_Change starting at row 2.
$r = 2
foreach ($row in $sheet2Data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $ws2.Cells.Item($r, 8).Value = $row[7]
    $ws2.Cells.Item($r, 9).Value = $row[8]
    $r++
}
